# Update TPM-derived ligand/receptor expression + specificity values
# (columns G-J, M-T) for Sema4c-Plxnb2 to reflect the new TPM normalization.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.20950266666667
$ws.Range("H2").Value = 39.628508
$ws.Range("I2").Value = 0.4617788189380864
$ws.Range("J2").Value = 0.4617788189380864
$ws.Range("M2").Value = 11.61289466666667
$ws.Range("N2").Value = 34.838684
$ws.Range("O2").Value = 0.09693042549509606
$ws.Range("P2").Value = 0.09693042549509606
$ws.Range("Q2").Value = 153.4005630670524
$ws.Range("R2").Value = 1380.605067603472
$ws.Range("S2").Value = 0.04476041740429164
$ws.Range("T2").Value = 0.04476041740429163
$ws.Range("G3").Value = 13.20950266666667
$ws.Range("H3").Value = 39.628508
$ws.Range("I3").Value = 0.4617788189380864
$ws.Range("J3").Value = 0.4617788189380864
$ws.Range("O3").Value = 0.2981108740043866
$ws.Range("P3").Value = 0.2981108740043866
$ws.Range("Q3").Value = 471.7855688253182
$ws.Range("R3").Value = 4246.070119427864
$ws.Range("S3").Value = 0.1376612873103463
$ws.Range("T3").Value = 0.1376612873103463
$ws.Range("G4").Value = 13.20950266666667
$ws.Range("H4").Value = 39.628508
$ws.Range("I4").Value = 0.4617788189380864
$ws.Range("J4").Value = 0.4617788189380864
$ws.Range("M4").Value = 27.39934733333333
$ws.Range("N4").Value = 82.198042
$ws.Range("O4").Value = 0.2286966748205465
$ws.Range("P4").Value = 0.2286966748205465
$ws.Range("Q4").Value = 361.9317516645929
$ws.Range("R4").Value = 3257.385764981336
$ws.Range("S4").Value = 0.1056072803936996
$ws.Range("T4").Value = 0.1056072803936996
$ws.Range("G5").Value = 13.20950266666667
$ws.Range("H5").Value = 39.628508
$ws.Range("I5").Value = 0.4617788189380864
$ws.Range("J5").Value = 0.4617788189380864
$ws.Range("M5").Value = 45.078635
$ws.Range("N5").Value = 135.235905
$ws.Range("O5").Value = 0.3762620256799708
$ws.Range("P5").Value = 0.3762620256799709
$ws.Range("Q5").Value = 595.4663492421932
$ws.Range("R5").Value = 5359.197143179739
$ws.Range("S5").Value = 0.1737498338297489
$ws.Range("T5").Value = 0.1737498338297489
$ws.Range("I6").Value = 0.2457011818953142
$ws.Range("J6").Value = 0.2457011818953142
$ws.Range("M6").Value = 11.61289466666667
$ws.Range("N6").Value = 34.838684
$ws.Range("O6").Value = 0.09693042549509606
$ws.Range("P6").Value = 0.09693042549509606
$ws.Range("Q6").Value = 81.62067661668758
$ws.Range("R6").Value = 734.586089550188
$ws.Range("S6").Value = 0.0238159201057608
$ws.Range("T6").Value = 0.0238159201057608
$ws.Range("I7").Value = 0.2457011818953142
$ws.Range("J7").Value = 0.2457011818953142
$ws.Range("O7").Value = 0.2981108740043866
$ws.Range("P7").Value = 0.2981108740043866
$ws.Range("S7").Value = 0.07324619407872288
$ws.Range("T7").Value = 0.07324619407872288
$ws.Range("I8").Value = 0.2457011818953142
$ws.Range("J8").Value = 0.2457011818953142
$ws.Range("M8").Value = 27.39934733333333
$ws.Range("N8").Value = 82.198042
$ws.Range("O8").Value = 0.2286966748205465
$ws.Range("P8").Value = 0.2286966748205465
$ws.Range("Q8").Value = 192.5750066967771
$ws.Range("R8").Value = 1733.175060270994
$ws.Range("S8").Value = 0.05619104329893663
$ws.Range("T8").Value = 0.05619104329893663
$ws.Range("I9").Value = 0.2457011818953142
$ws.Range("J9").Value = 0.2457011818953142
$ws.Range("M9").Value = 45.078635
$ws.Range("N9").Value = 135.235905
$ws.Range("O9").Value = 0.3762620256799708
$ws.Range("P9").Value = 0.3762620256799709
$ws.Range("Q9").Value = 316.8330373492317
$ws.Range("R9").Value = 2851.497336143085
$ws.Range("S9").Value = 0.09244802441189388
$ws.Range("T9").Value = 0.0924480244118939
$ws.Range("G10").Value = 7.459692333333334
$ws.Range("H10").Value = 22.379077
$ws.Range("I10").Value = 0.2607765032684172
$ws.Range("J10").Value = 0.2607765032684171
$ws.Range("M10").Value = 11.61289466666667
$ws.Range("N10").Value = 34.838684
$ws.Range("O10").Value = 0.09693042549509606
$ws.Range("P10").Value = 0.09693042549509606
$ws.Range("Q10").Value = 86.6286213127409
$ws.Range("R10").Value = 779.6575918146681
$ws.Range("S10").Value = 0.02527717742093099
$ws.Range("T10").Value = 0.02527717742093098
$ws.Range("G11").Value = 7.459692333333334
$ws.Range("H11").Value = 22.379077
$ws.Range("I11").Value = 0.2607765032684172
$ws.Range("J11").Value = 0.2607765032684171
$ws.Range("O11").Value = 0.2981108740043866
$ws.Range("P11").Value = 0.2981108740043866
$ws.Range("Q11").Value = 266.4275317211185
$ws.Range("R11").Value = 2397.847785490067
$ws.Range("S11").Value = 0.07774031130915561
$ws.Range("T11").Value = 0.0777403113091556
$ws.Range("G12").Value = 7.459692333333334
$ws.Range("H12").Value = 22.379077
$ws.Range("I12").Value = 0.2607765032684172
$ws.Range("J12").Value = 0.2607765032684171
$ws.Range("M12").Value = 27.39934733333333
$ws.Range("N12").Value = 82.198042
$ws.Range("O12").Value = 0.2286966748205465
$ws.Range("P12").Value = 0.2286966748205465
$ws.Range("Q12").Value = 204.3907012408038
$ws.Range("R12").Value = 1839.516311167234
$ws.Range("S12").Value = 0.05963871916881639
$ws.Range("T12").Value = 0.05963871916881638
$ws.Range("G13").Value = 7.459692333333334
$ws.Range("H13").Value = 22.379077
$ws.Range("I13").Value = 0.2607765032684172
$ws.Range("J13").Value = 0.2607765032684171
$ws.Range("M13").Value = 45.078635
$ws.Range("N13").Value = 135.235905
$ws.Range("O13").Value = 0.3762620256799708
$ws.Range("P13").Value = 0.3762620256799709
$ws.Range("Q13").Value = 336.2727479066317
$ws.Range("R13").Value = 3026.454731159685
$ws.Range("S13").Value = 0.09812029536951417
$ws.Range("T13").Value = 0.09812029536951417
$ws.Range("G14").Value = 0.9080446666666666
$ws.Range("H14").Value = 2.724134
$ws.Range("I14").Value = 0.03174349589818232
$ws.Range("J14").Value = 0.03174349589818231
$ws.Range("M14").Value = 11.61289466666667
$ws.Range("N14").Value = 34.838684
$ws.Range("O14").Value = 0.09693042549509606
$ws.Range("P14").Value = 0.09693042549509606
$ws.Range("Q14").Value = 10.54502706662844
$ws.Range("R14").Value = 94.90524359965599
$ws.Range("S14").Value = 0.003076910564112648
$ws.Range("T14").Value = 0.003076910564112648
$ws.Range("G15").Value = 0.9080446666666666
$ws.Range("H15").Value = 2.724134
$ws.Range("I15").Value = 0.03174349589818232
$ws.Range("J15").Value = 0.03174349589818231
$ws.Range("O15").Value = 0.2981108740043866
$ws.Range("P15").Value = 0.2981108740043866
$ws.Range("Q15").Value = 32.43137765233022
$ws.Range("R15").Value = 291.882398870972
$ws.Range("S15").Value = 0.00946308130616179
$ws.Range("T15").Value = 0.009463081306161789
$ws.Range("G16").Value = 0.9080446666666666
$ws.Range("H16").Value = 2.724134
$ws.Range("I16").Value = 0.03174349589818232
$ws.Range("J16").Value = 0.03174349589818231
$ws.Range("M16").Value = 27.39934733333333
$ws.Range("N16").Value = 82.198042
$ws.Range("O16").Value = 0.2286966748205465
$ws.Range("P16").Value = 0.2286966748205465
$ws.Range("Q16").Value = 24.87983121618089
$ws.Range("R16").Value = 223.918480945628
$ws.Range("S16").Value = 0.007259631959093954
$ws.Range("T16").Value = 0.007259631959093953
$ws.Range("G17").Value = 0.9080446666666666
$ws.Range("H17").Value = 2.724134
$ws.Range("I17").Value = 0.03174349589818232
$ws.Range("J17").Value = 0.03174349589818231
$ws.Range("M17").Value = 45.078635
$ws.Range("N17").Value = 135.235905
$ws.Range("O17").Value = 0.3762620256799708
$ws.Range("P17").Value = 0.3762620256799709
$ws.Range("Q17").Value = 40.93341409236333
$ws.Range("R17").Value = 368.40072683127
$ws.Range("S17").Value = 0.01194387206881392
$ws.Range("T17").Value = 0.01194387206881392
